$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row ---
# J1 becomes "Giornata 9" (6 new Giornata headers shift the averages column to P)
$ws.Range("J1").Value = "Giornata 9"
$ws.Range("K1").Value = "Giornata 10"
$ws.Range("L1").Value = "Giornata 11"
$ws.Range("M1").Value = "Giornata 12"
$ws.Range("N1").Value = "Giornata 13"
$ws.Range("O1").Value = "Giornata 14"
$ws.Range("P1").Value = "Total average "

# Match the header style (bold, centered, thin border) used by A1:I1
$hdr = $ws.Range("J1:P1")
$hdr.Font.Bold = $true
$hdr.HorizontalAlignment = -4108
$hdr.VerticalAlignment = -4160
$hdr.Borders.LineStyle = 1

# --- Data rows: Giornata 9 (col J) moves in, Giornata 10-14 (K-O) are new, P holds the updated overall average ---
# Row 2: ATALANTA
$ws.Range("J2").Value = 5.633333333333334
$ws.Range("L2").Value = 6.615384615384615
$ws.Range("M2").Value = 6.428571428571429
$ws.Range("N2").Value = 6.714285714285714
$ws.Range("O2").Value = 6.178571428571429
$ws.Range("P2").Value = 6.106146097492251

# Row 3: BENEVENTO
$ws.Range("J3").Value = 6.15625
$ws.Range("K3").Value = 5.75
$ws.Range("L3").Value = 6
$ws.Range("M3").Value = 6.071428571428571
$ws.Range("N3").Value = 6.178571428571429
$ws.Range("O3").Value = 6.46875
$ws.Range("P3").Value = 5.971989468864469

# Row 4: BOLOGNA
$ws.Range("J4").Value = 6.153846153846154
$ws.Range("K4").Value = 5.615384615384615
$ws.Range("L4").Value = 4.75
$ws.Range("M4").Value = 5.875
$ws.Range("N4").Value = 5.807692307692307
$ws.Range("O4").Value = 6.071428571428571
$ws.Range("P4").Value = 5.895532443746729

# Row 5: CAGLIARI
$ws.Range("J5").Value = 5.961538461538462
$ws.Range("K5").Value = 6.045454545454546
$ws.Range("L5").Value = 5.678571428571429
$ws.Range("M5").Value = 5.961538461538462
$ws.Range("N5").Value = 6
$ws.Range("O5").Value = 5.923076923076923
$ws.Range("P5").Value = 5.96059892488464

# Row 6: CROTONE
$ws.Range("J6").Value = 5.866666666666666
$ws.Range("K6").Value = 5.291666666666667
$ws.Range("L6").Value = 6.458333333333333
$ws.Range("M6").Value = 5.833333333333333
$ws.Range("N6").Value = 5.730769230769231
$ws.Range("O6").Value = 6.269230769230769
$ws.Range("P6").Value = 5.849058084772371

# Row 7: FIORENTINA
$ws.Range("J7").Value = 5.5
$ws.Range("K7").Value = 5.692307692307693
$ws.Range("L7").Value = 5.5
$ws.Range("M7").Value = 6
$ws.Range("N7").Value = 5.964285714285714
$ws.Range("O7").Value = 7.038461538461538
$ws.Range("P7").Value = 5.883710931925217

# Row 8: GENOA
$ws.Range("J8").Value = 5.615384615384615
$ws.Range("K8").Value = 5.84375
$ws.Range("L8").Value = 5.642857142857143
$ws.Range("M8").Value = 6.142857142857143
$ws.Range("N8").Value = 5.5
$ws.Range("O8").Value = 6.3
$ws.Range("P8").Value = 5.801914900575615

# Row 9: INTER
$ws.Range("J9").Value = 6.772727272727272
$ws.Range("K9").Value = 6.428571428571429
$ws.Range("L9").Value = 6.03125
$ws.Range("M9").Value = 6.076923076923077
$ws.Range("N9").Value = 6.076923076923077
$ws.Range("O9").Value = 6.041666666666667
$ws.Range("P9").Value = 6.070251693334678

# Row 10: JUVENTUS
$ws.Range("J10").Value = 5.535714285714286
$ws.Range("K10").Value = 5.892857142857143
$ws.Range("L10").Value = 6.291666666666667
$ws.Range("M10").Value = 6.153846153846154
$ws.Range("N10").Value = 6.892857142857143
$ws.Range("O10").Value = 5.25
$ws.Range("P10").Value = 6.101627539127541

# Row 11: LAZIO
$ws.Range("J11").Value = 5
$ws.Range("K11").Value = 6.133333333333334
$ws.Range("L11").Value = 5.59375
$ws.Range("M11").Value = 5.833333333333333
$ws.Range("N11").Value = 6.375
$ws.Range("O11").Value = 6
$ws.Range("P11").Value = 5.921322278911563

# Row 12: MILAN
$ws.Range("J12").Value = 6.5
$ws.Range("K12").Value = 6.192307692307693
$ws.Range("L12").Value = 5.807692307692307
$ws.Range("M12").Value = 5.807692307692307
$ws.Range("N12").Value = 6.607142857142857
$ws.Range("O12").Value = 6.333333333333333
$ws.Range("P12").Value = 6.270329670329671

# Row 13: NAPOLI
$ws.Range("J13").Value = 6.791666666666667
$ws.Range("K13").Value = 6.541666666666667
$ws.Range("L13").Value = 6.033333333333333
$ws.Range("M13").Value = 5.884615384615385
$ws.Range("N13").Value = 5.133333333333334
$ws.Range("O13").Value = 5.5625
$ws.Range("P13").Value = 6.114089779268352

# Row 14: PARMA
$ws.Range("J14").Value = 6.133333333333334
$ws.Range("K14").Value = 5.5
$ws.Range("L14").Value = 6.115384615384615
$ws.Range("M14").Value = 5.875
$ws.Range("N14").Value = 5.1
$ws.Range("O14").Value = 5.6
$ws.Range("P14").Value = 5.752191261119831

# Row 15: ROMA
$ws.Range("J15").Value = 4.964285714285714
$ws.Range("K15").Value = 5.909090909090909
$ws.Range("L15").Value = 6.615384615384615
$ws.Range("M15").Value = 6.0625
$ws.Range("N15").Value = 5.416666666666667
$ws.Range("O15").Value = 6.125
$ws.Range("P15").Value = 6.122611019932449

# Row 16: SAMPDORIA
$ws.Range("J16").Value = 5.966666666666667
$ws.Range("K16").Value = 5.8
$ws.Range("L16").Value = 5.833333333333333
$ws.Range("M16").Value = 6.269230769230769
$ws.Range("N16").Value = 6.576923076923077
$ws.Range("O16").Value = 5.8
$ws.Range("P16").Value = 6.022928113553113

# Row 17: SASSUOLO
$ws.Range("J17").Value = 5.375
$ws.Range("K17").Value = 5.857142857142857
$ws.Range("L17").Value = 5.857142857142857
$ws.Range("M17").Value = 6.269230769230769
$ws.Range("N17").Value = 5.8
$ws.Range("O17").Value = 6.214285714285714
$ws.Range("P17").Value = 6.01854722658294

# Row 18: SPEZIA
$ws.Range("J18").Value = 6.266666666666667
$ws.Range("K18").Value = 5.821428571428571
$ws.Range("L18").Value = 5.5
$ws.Range("M18").Value = 5.964285714285714
$ws.Range("N18").Value = 5.964285714285714
$ws.Range("O18").Value = 5.571428571428571
$ws.Range("P18").Value = 5.882568027210882

# Row 19: TORINO
$ws.Range("J19").Value = 6.038461538461538
$ws.Range("K19").Value = 5.791666666666667
$ws.Range("L19").Value = 5.2
$ws.Range("M19").Value = 5.833333333333333
$ws.Range("N19").Value = 5.769230769230769
$ws.Range("O19").Value = 6.409090909090909
$ws.Range("P19").Value = 5.813456781313924

# Row 20: UDINESE
$ws.Range("J20").Value = 6.642857142857143
$ws.Range("L20").Value = 6.321428571428571
$ws.Range("M20").Value = 6.125
$ws.Range("N20").Value = 6.357142857142857
$ws.Range("O20").Value = 5.464285714285714
$ws.Range("P20").Value = 6.047231614539307

# Row 21: VERONA
$ws.Range("J21").Value = 6.53125
$ws.Range("K21").Value = 5.8
$ws.Range("L21").Value = 6.307692307692307
$ws.Range("M21").Value = 5.566666666666666
$ws.Range("N21").Value = 5.933333333333334
$ws.Range("O21").Value = 5.78125
$ws.Range("P21").Value = 6.062921899529043
